$d = $word.ActiveDocument
$wNs = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'

# ------------------------------------------------------------------
# 1) "Must change goldenrod city flower shop ..." paragraph: split out
#    "sudowoodo" into its own run wrapped in proofErr spellStart/spellEnd.
# ------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2)
$xml2 = '<w:p xmlns:w="' + $wNs + '" w:rsidR="006F1D20" w:rsidRDefault="00996BEE" w:rsidP="00996BEE">' +
        '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
        '<w:r><w:t xml:space="preserve">Must change goldenrod city flower shop to REMOVE the check on beating </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>sudowoodo</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t>, or the location will lock itself</w:t></w:r>' +
        '</w:p>'
$p2.Range.InsertXML($xml2)

# ------------------------------------------------------------------
# 2) "Need to actually figure out ..." paragraph: split out "despawn"
#    and "beated" into their own runs wrapped in proofErr markers.
# ------------------------------------------------------------------
$p4 = $d.Paragraphs.Item(4)
$xml4 = '<w:p xmlns:w="' + $wNs + '" w:rsidR="00A641A6" w:rsidRDefault="00A641A6" w:rsidP="00996BEE">' +
        '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
        '<w:r><w:t xml:space="preserve">Need to actually figure out how to make it so the director in the underground warehouse won' + [char]8217 + 't </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>despawn</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> when team rocket is </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>beated</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t>(and also make sure he still spawns in on the rocket trigger)</w:t></w:r>' +
        '</w:p>'
$p4.Range.InsertXML($xml4)

# ------------------------------------------------------------------
# 3) Remove the "Modify bug catching contest ..." paragraph entirely
#    (paragraph 7), and strip the "Don't let the shop in goldenrod
#    sell rock smash ..." paragraph (old paragraph 8) down to just the
#    bookmark it carries, then drop the now-redundant trailing empty
#    paragraph by merging it away.
# ------------------------------------------------------------------
$p7 = $d.Paragraphs.Item(7)
$r7 = $d.Range($p7.Range.Start, $p7.Range.End)
$r7.Delete()

# After that deletion, "Don't let the shop ..." is paragraph 7 and the
# trailing empty paragraph is paragraph 8. Merge the trailing empty
# paragraph into the bookmark paragraph by deleting the paragraph mark
# between them.
$p7 = $d.Paragraphs.Item(7)
$p8 = $d.Paragraphs.Item(8)
$rMerge = $d.Range($p7.Range.End - 1, $p8.Range.End)
$rMerge.Delete()

# Replace the bookmark paragraph's content with just the bookmark.
$p7 = $d.Paragraphs.Item(7)
$xml7 = '<w:p xmlns:w="' + $wNs + '" w:rsidR="0071690E" w:rsidRDefault="0071690E" w:rsidP="00996BEE">' +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
        '<w:bookmarkEnd w:id="0"/>' +
        '</w:p>'
$p7.Range.InsertXML($xml7)

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
